# Add the latest dividend entry (17/09/2025, Gross Dividend 0.005) to the
# top of the DividendHistory table, pushing all existing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new, most-recent dividend record right under the header.
$ws.Rows.Item(2).Insert()

# XD Date / Pay Date are both 17/09/2025 for the new entry.
$ws.Range("A2").Value = "17/09/2025"
$ws.Range("B2").Value = "17/09/2025"

# Gross Dividend stays the same as the previous period (0.005). Copy it from
# the row below (already stored as text) instead of retyping the literal, so
# the numeric-looking string isn't auto-converted into a Number cell.
$ws.Range("C3").Copy()
$ws.Range("C2").PasteSpecial()
